$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("optimization_parameters")

# Row 1: drop the redundant trailing "value" header cells (C1:F1)
$ws.Range("C1:F1").Clear()

# Rename the "Model" parameter row to "production_function" and give its
# label cell the same bold header style as A1/B1.
$ws.Range("A1").Copy()
$ws.Range("A8").PasteSpecial(-4122)
$ws.Range("A8").Value = "production_function"

# Insert a new row for the "L_curve" parameter right after it.
$ws.Rows.Item(9).Insert()
$ws.Range("A1").Copy()
$ws.Range("A9").PasteSpecial(-4122)
$ws.Range("A9").Value = "L_curve"
$ws.Range("B2").Copy()
$ws.Range("B9").PasteSpecial(-4122)
$ws.Range("B9").Value = 1

# Remove the obsolete "Deletion" row (now pushed down to row 17).
$ws.Rows.Item(17).Delete()

# Make this the active sheet / tab, with B25 selected, matching the new
# workbook view state.
$ws.Activate()
$ws.Range("B25").Select()
